$excel.DisplayAlerts = $false | Out-Null

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CFG 3 Frasi")

# --- Clear existing content on "CFG 3 Frasi" so we can rebuild it cleanly ---
$ws.Cells.Clear() | Out-Null

# --- Header row ---
$ws.Cells.Item(1,1).Value = "Regole"
$ws.Cells.Item(1,2).Value = "Conta se usate"
$ws.Cells.Item(1,3).Value = "in frasi"
$ws.Cells.Item(1,4).Value = "in frasi2"
$ws.Cells.Item(1,5).Value = "Status"
$ws.Cells.Item(1,7).Value = "ID"
$ws.Cells.Item(1,8).Value = "Frasi"

# --- Grammar-rule rows (A..E), rows 2-15 ---
$rules = @(
    @("S -> NP VP",        0,   $null, $null, $null),
    @("S -> NP VP PP",     1,   3,     $null, "ok"),
    @("S -> VP NP",        1,   1,     $null, "ok"),
    @("S -> VP NP PP",     1,   2,     $null, "ok"),
    @("NP -> NN",          1,   3,     $null, "ok"),
    @("NP -> NNS",         1,   1,     $null, "ok"),
    @("NP -> DT NN",       1,   2,     $null, "ok"),
    @("NP -> PRP`$ JJ NP", 1,   3,     $null, "ok"),
    @("VP -> PRP VBP VBG", 1,   1,     $null, "ok"),
    @("VP -> EX VBZ",      1,   2,     $null, "ok"),
    @("VP -> VBZ VBG",     1,   3,     $null, "ok"),
    @("PP -> IN PP",       2,   2,     3,     "ambiguità"),
    @("PP -> PRP`$ NN",    1,   2,     $null, "ok"),
    @("PP -> IN RB",       1,   3,     $null, "ok")
)

$r = 2
foreach ($rule in $rules) {
    $ws.Cells.Item($r,1).Value = $rule[0]
    $ws.Cells.Item($r,2).Value = $rule[1]
    if ($rule[2] -ne $null) { $ws.Cells.Item($r,3).Value = $rule[2] }
    if ($rule[3] -ne $null) { $ws.Cells.Item($r,4).Value = $rule[3] }
    if ($rule[4] -ne $null) { $ws.Cells.Item($r,5).Value = $rule[4] }
    $r = $r + 1
}

# --- Sentences moved over from the old "Frasi" sheet (G..H), rows 2-13 ---
$sentences = @(
    "You are imagining things.",
    "There is a price on my head",
    "Your big opportinity is flying out of here",
    "You are eating too much",
    "Lorenzo is sleeping too much",
    "She is running at home",
    "It's raining over my head",
    "A piano is falling over my car",
    "A brick has fallen beside my foot",
    "Your girlfriend is waiting out of here",
    "Your best friend is playing in the garden",
    "Angus is giving the dog a bone"
)

$r = 2
foreach ($s in $sentences) {
    $ws.Cells.Item($r,7).Value = $r - 1
    $ws.Cells.Item($r,8).Value = $s
    $r = $r + 1
}

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(7).ColumnWidth = 10.166666666666666
$ws.Columns.Item(8).ColumnWidth = 33.498697916666664

# --- Tables ---
$tbl1 = $ws.ListObjects.Add(1, $ws.Range("A1:E15"), 0, 1)
$tbl1.Name = "CFG3Frasi"
$tbl1.TableStyle = "TableStyleLight8"

$tbl2 = $ws.ListObjects.Add(1, $ws.Range("G1:H13"), 0, 1)
$tbl2.Name = "Frasi"
$tbl2.TableStyle = "TableStyleLight8"

# --- Remove the now-redundant standalone "Frasi" sheet ---
$wb.Worksheets.Item("Frasi").Delete() | Out-Null

# --- Make "CFG 3 Frasi" the active sheet/tab (re-fetch: the old handle is stale after the delete) ---
$ws = $wb.Worksheets.Item("CFG 3 Frasi")
$ws.Activate() | Out-Null
$ws.Range("H21").Select() | Out-Null
